$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 226, shifting existing rows 226:316 down to 227:317
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with a new weekly price record
$ws.Range("A226").Value = 10
$ws.Range("B226").Value = "Vega Modelo de Temuco"
$ws.Range("C226").Value = "La Araucanía"
$ws.Range("D226").Value = 44924
$ws.Range("E226").Value = 9
$ws.Range("F226").Value = 100112039
$ws.Range("G226").Value = "Ciboulette"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 65
$ws.Range("K226").Value = 5000
$ws.Range("L226").Value = 5000
$ws.Range("M226").Value = 5000
$ws.Range("N226").Value = "$/docena de atados"
$ws.Range("O226").Value = "Provincia de Cautín"
$ws.Range("P226").Value = 1667
$ws.Range("Q226").Value = 3
$ws.Range("R226").Value = "Hortaliza"
